$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 55: give C55 a (blank) percent-formatted cell, like the other rows in this table
$ws.Range("C55").NumberFormat = "0%"

# New row 58: "Agregar los comandos abajo, en la pantalla (teclas rapidas)" / "Lucas"
$ws.Range("A58").Value = "Agregar los comandos abajo, en la pantalla (teclas rapidas)"
$ws.Range("B58").Value = "Lucas"

# Row 59 gets replaced with a new task: "Agregar signo $ en los totales" (no responsible)
$ws.Range("A59").Value = "Agregar signo `$ en los totales"

# New row 60: "En reportes mensual y anual poner mes y año, no dia!" / "Agustina"
$ws.Range("A60").Value = "En reportes mensual y anual poner mes y año, no dia!"
$ws.Range("B60").Value = "Agustina"

# The old row 59 content ("Ivan: preguntar reportes...") moves down to row 66
$ws.Range("A66").Value = "Ivan: preguntar reportes - preguntar autorizacion requerida en que funciones - preguntar login"

# Update the view: scrolled/selected position moved down a bit
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("B59").Select()
